$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DashboardButton")

# Data for the new rows (row 4 through row 11), columns A..I:
# A = Editorial Role, B = PostTilte, C = Logged_in author, D = Authorname,
# E = Post_type, F = type, G = Expected_Buttons, H = Actual_buttons, I = Status
# A blank string value ("") below stands for an empty cell that must still be
# written out as a (shared-)string cell, matching the source row 3 behaviour.

$rows = @(
    @("Coordinator", "Future brand club collaborator role isha", "unbrandcord", "col", "schedule", "Club", "", "", "Pass"),
    @("Coordinator", "Future brand club coordinator role isha", "unbrandcord", "shobha", "schedule", "Club", "", "", "Pass"),
    @("Coordinator", "Future brand longform coordinator role isha", "unbrandcord", "shobha", "schedule", "Club", "", "", "Pass"),
    @("Coordinator", "Future slideshow post editor role isha", "unbrandcord", "shobhaeditor", "schedule", "normal", "Repost,Editar,Pasar a borrador,Destacar", "Repost,Editar,Pasar a borrador,Destacar", "Pass"),
    @("Coordinator", "coordinator post longform :- Sumit", "unbrandcord", "coordinatorIsha", "schedule", "normal", "Repost,Editar,Pasar a borrador,Destacar", "Repost,Editar,Pasar a borrador,Destacar", "Pass"),
    @("Director", "Future normal longform editor role isha", "director", "shobhaeditor", "schedule", "normal", "Repost,Editar,Pasar a borrador,Destacar", "Repost,Editar,Pasar a borrador,Destacar", "Pass"),
    @("Director", "Future brand longform collaborator role isha", "director", "col", "schedule", "Club", "", "", "Pass"),
    @("Director", "Future brand club collaborator role isha", "director", "col", "schedule", "Club", "", "", "Pass")
)

$startRow = 4
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $val = $rowData[$c - 1]
        $cell = $ws.Cells.Item($r, $c)
        if ($val -eq "") {
            # Writing a plain "" leaves the cell completely blank/absent, so
            # use Excel's text-prefix escape to force an empty text cell,
            # then reset the style so no quote-prefix formatting sticks.
            $cell.Value = "'"
            $cell.Style = "Normal"
        } else {
            $cell.Value = $val
        }
    }
}
